$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry represents a data row whose "剩余" (E column, remaining days)
# and "开始时间" (F column, start date) values need to be updated to reflect
# one additional day having passed (auto-update run).
$updates = @(
    @{Row=2; E=11; F=20251229},
    @{Row=3; E=11; F=20251229},
    @{Row=4; E=11; F=20251229},
    @{Row=5; E=3; F=20251225},
    @{Row=6; E=11; F=20251229},
    @{Row=7; E=3; F=20251225},
    @{Row=8; E=11; F=20251229},
    @{Row=9; E=3; F=20251225},
    @{Row=10; E=4; F=20251229},
    @{Row=11; E=11; F=20251229},
    @{Row=12; E=3; F=20251225},
    @{Row=13; E=11; F=20251229},
    @{Row=14; E=11; F=20251229},
    @{Row=15; E=11; F=20251229},
    @{Row=16; E=7; F=20251229},
    @{Row=17; E=3; F=20251225},
    @{Row=18; E=6; F=20251228},
    @{Row=19; E=6; F=20251228},
    @{Row=20; E=6; F=20251228},
    @{Row=21; E=6; F=20251228},
    @{Row=22; E=3; F=20251225},
    @{Row=23; E=3; F=20251225},
    @{Row=24; E=3; F=20251225},
    @{Row=25; E=3; F=20251225},
    @{Row=26; E=3; F=20251225},
    @{Row=27; E=5; F=20251230},
    @{Row=28; E=6; F=20251228},
    @{Row=29; E=6; F=20251228},
    @{Row=30; E=6; F=20251228},
    @{Row=31; E=6; F=20251228},
    @{Row=32; E=6; F=20251228},
    @{Row=33; E=6; F=20251228},
    @{Row=34; E=6; F=20251228},
    @{Row=35; E=6; F=20251228},
    @{Row=37; E=6; F=20251228},
    @{Row=38; E=6; F=20251228},
    @{Row=39; E=6; F=20251228},
    @{Row=40; E=4; F=20251229},
    @{Row=41; E=4; F=20251229},
    @{Row=42; E=6; F=20251228},
    @{Row=43; E=3; F=20251225},
    @{Row=44; E=4; F=20251229},
    @{Row=45; E=3; F=20251225},
    @{Row=46; E=4; F=20251229},
    @{Row=47; E=6; F=20251228},
    @{Row=48; E=4; F=20251229},
    @{Row=49; E=5; F=20251230},
    @{Row=50; E=1; F=20251223},
    @{Row=51; E=1; F=20251223},
    @{Row=52; E=1; F=20251223},
    @{Row=53; E=1; F=20251223},
    @{Row=54; E=1; F=20251223},
    @{Row=55; E=1; F=20251223},
    @{Row=56; E=1; F=20251223},
    @{Row=57; E=1; F=20251223},
    @{Row=58; E=5; F=20251227},
    @{Row=59; E=5; F=20251227},
    @{Row=60; E=5; F=20251227},
    @{Row=61; E=5; F=20251230},
    @{Row=62; E=5; F=20251227},
    @{Row=63; E=5; F=20251227},
    @{Row=64; E=5; F=20251227},
    @{Row=65; E=6; F=20251228},
    @{Row=66; E=6; F=20251228},
    @{Row=67; E=6; F=20251228},
    @{Row=68; E=6; F=20251228},
    @{Row=69; E=6; F=20251228},
    @{Row=70; E=7; F=20251229},
    @{Row=71; E=7; F=20251229},
    @{Row=72; E=7; F=20251229},
    @{Row=73; E=7; F=20251229},
    @{Row=74; E=7; F=20251229},
    @{Row=75; E=7; F=20251229},
    @{Row=76; E=7; F=20251229},
    @{Row=77; E=10; F=20260101},
    @{Row=78; E=10; F=20260101},
    @{Row=79; E=10; F=20260101},
    @{Row=80; E=10; F=20260101},
    @{Row=81; E=10; F=20260101},
    @{Row=82; E=10; F=20260101},
    @{Row=83; E=10; F=20260101},
    @{Row=84; E=10; F=20260101},
    @{Row=85; E=10; F=20260101},
    @{Row=86; E=10; F=20260101},
    @{Row=87; E=4; F=20251229},
    @{Row=88; E=4; F=20251229},
    @{Row=89; E=4; F=20251229},
    @{Row=90; E=4; F=20251229},
    @{Row=91; E=3; F=20251225},
    @{Row=92; E=4; F=20251229},
    @{Row=93; E=10; F=20260101},
    @{Row=94; E=7; F=20260101},
    @{Row=95; E=9; F=20251231},
    @{Row=96; E=7; F=20251229},
    @{Row=97; E=7; F=20251229},
    @{Row=98; E=7; F=20251229},
    @{Row=99; E=7; F=20251229}
)

foreach ($u in $updates) {
    $row = $u.Row
    $ws.Cells.Item($row, 5).Value = $u.E
    $ws.Cells.Item($row, 6).Value = $u.F
}
